$d = $word.ActiveDocument

# 1) Ativação date update.
# This run is immediately followed (after a <w:br/>) by the "Departamento: ..."
# run which has identical (empty) run formatting. A plain text replacement
# here would cause the engine to coalesce the two runs together (merging the
# <w:br/> into the middle of a single run), which does not match the
# original run layout. To keep the two runs distinct (as in the target),
# we briefly touch the formatting of the following run so it is treated as
# a separate run, then restore it.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $pp = $paras.Item($i)
    if ($pp.Range.Text.StartsWith("Créditos-aula")) {
        $full = $pp.Range.Text
        $idx = $full.IndexOf("Ativação: 01/01/2012")
        $start = $pp.Range.Start + $idx
        $end = $start + "Ativação: 01/01/2012".Length
        $rngDate = $d.Range($start, $end)
        $rngDate.Text = "Ativação: 01/01/2021"

        $full2 = $pp.Range.Text
        $idx2 = $full2.IndexOf("Departamento")
        $start2 = $pp.Range.Start + $idx2
        $end2 = $start2 + "Departamento: Engenharia Química".Length
        $rngDept = $d.Range($start2, $end2)
        $rngDept.Font.Bold = 1
        $rngDept.Font.Bold = 0
        break
    }
}

# 2) Objetivos paragraph text
$d.Content.Find.Execute(
    "Apresentar os conceitos básicos da Logística Empresarial e da Gestão da Cadeia de Suprimentos. Capacitar o aluno para aplicação de técnicas e métodos quantitativos para otimização dos problemas em Logística e Cadeias de Suprimentos.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Apresentar os conceitos de Logística, Logística Reversa e Gestão da Cadeia de Suprimentos. Capacitar o aluno para aplicação de técnicas e métodos quantitativos para otimização dos problemas em Logística e Cadeias de Suprimentos.",
    2) | Out-Null

# 3) Programa resumido paragraph text
$d.Content.Find.Execute(
    "Logística Integrada, Serviço ao Cliente, Administração do Transporte, Custos Logísticos, Armazenagem e Localização das Instalações, Tecnologia de Informação Aplicada à Logística, Canais de Distribuição e Distribuição Física, Roterização de Veículos e Operadores Logísticos, Organização e Controle de Estoques.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "1. Introdução: 2. Gestão estratégica3. Gestão dos relacionamentos4. Gestão global de suprimentos5. Avaliação de desempenho6. Mapeamento e análise de processos7. Gestão de demanda8. Gestão e coordenação de estoques9. Gestão da logística10. Logística reversa",
    2) | Out-Null

# 4) Programa paragraph: collapse the many runs/breaks into a single run with
#    the new plain text (no line breaks). Find it by locating the paragraph
#    that begins with "1. Logística Integrada".
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("1. Logística Integrada")) {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "1. Introdução: A concorrência entre cadeias de suprimento. Definição operacional. A globalização e a gestão de cadeia de suprimentos. Governança das cadeias de suprimentos2. Gestão estratégica: Estratégia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e híbridos. Custo de transação e a decisão estratégica de comprar ou fazer. Padronização. Integração de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gestão dos relacionamentos: Confiança entre parceiros. Negociação. Gestão do relacionamento com clientes. Segmentação de produtos. Gestão do relacionamento com fornecedores4. Gestão global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopetição. Ética e responsabilidade social na gestão global de suprimentos5. Avaliação de desempenho: O que é medição de desempenho? Porque medir desempenho. Características de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e análise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). Análise e melhoramento de processos.7. Gestão de demanda: Ações sobre a demanda para redução de variabilidade. Causas da variabilidade da demanda. Previsão de demanda. Processo de previsão de vendas. Métodos usados em previsões. Método Delphi. Incerteza de previsão8. Gestão e coordenação de estoques: Definição de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consignação9. Gestão da logística: Centralização versus descentralização. Pontos de armazenagem/distribuição. Funções dos armazéns. Sistemas logísticos escalonados. Localização de unidades logísticas. Gestão de transportes na cadeia de suprimentos.10. Logística reversa: Conceito, importância, estrutura e tendências. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motivação empresarial. Gerenciamento integrado de resíduos. Modelos de roteirização. Programação de frotas de veículos."
        break
    }
}

# 5) Critério text
$d.Content.Find.Execute(
    "Média aritmética de duas provas teóricas.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Média de Provas e trabalhos (MF).", 2) | Out-Null

# 6) Norma de recuperação text
$d.Content.Find.Execute(
    "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR",
    2) | Out-Null

# 7) Bibliografia paragraph: collapse the many runs/breaks into a single run
#    with the new plain text (no line breaks). Find it by locating the
#    paragraph that begins with "CORRÊA, H. L.".
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("CORRÊA, H. L.")) {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "CORRÊA, HENRIQUE LUIZ. Gestão de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administração de cadeias de suprimento e logística: o essencial. Editora Atlas 2014PIRES, SÉRGIO. Gestão da cadeia de suprimentos (Supply Chain Management): conceitos, estratégias, práticas e casos. Editora Atlas segunda edição. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gestão da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012"
        break
    }
}
